$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.791.11'
$ws.Range('E2').Value = '  -0.14%  '
$ws.Range('D3').Value = '2.384.34'
$ws.Range('E3').Value = '  -1.23%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '555.33'
$ws.Range('E5').Value = '  +0.69%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '133.39'
$ws.Range('E6').Value = '  -2.64%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('E8').Value = '  -0.74%  '
$ws.Range('E9').Value = '  -0.27%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.61'
$ws.Range('E10').Value = '  -1.61%  '
$ws.Range('E11').Value = '  +1.16%  '
$ws.Range('E12').Value = '  -2.80%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '24.40'
$ws.Range('E13').Value = '  -4.65%  '
$ws.Range('D14').Value = '2.803.52'
$ws.Range('E14').Value = '  -1.40%  '
$ws.Range('D15').Value = '59.721.59'
$ws.Range('E15').Value = '  -0.15%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000137'
$ws.Range('E16').Value = '  -0.52%  '
$ws.Range('D17').Value = '2.363.77'
$ws.Range('E17').Value = '  -1.48%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.11'
$ws.Range('E18').Value = '  -1.80%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.48'
$ws.Range('E19').Value = '  +1.56%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '319.75'
$ws.Range('E20').Value = '  -2.66%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.70'
$ws.Range('E21').Value = '  +0.73%  '
$ws.Range('E22').Value = '  +0.09%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '64.14'
$ws.Range('E23').Value = '  -3.44%  '
$ws.Range('E24').Value = '  +0.76%  '
$ws.Range('E25').Value = '  -0.05%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.39'
$ws.Range('E26').Value = '  -2.87%  '
$ws.Range('E27').Value = '  -0.10%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.80'
$ws.Range('E28').Value = '  +1.82%  '
$ws.Range('E29').Value = '  -2.11%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '169.62'
$ws.Range('E30').Value = '  +0.90%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.05'
$ws.Range('E31').Value = '  -0.88%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.09'
$ws.Range('E32').Value = '  +7.86%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.395'
$ws.Range('E33').Value = '  -3.19%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '18.18'
$ws.Range('E34').Value = '  -2.31%  '
$ws.Range('E36').Value = '  +1.31%  '
$ws.Range('E37').Value = '  +0.05%  '
$ws.Range('E38').Value = '  -1.87%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '318.56'
$ws.Range('E39').Value = '  +1.52%  '
$ws.Range('E40').Value = '  -1.80%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '38.62'
$ws.Range('E41').Value = '  -2.43%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '145.72'
$ws.Range('E42').Value = '  +5.06%  '
$ws.Range('E43').Value = '  -4.23%  '
$ws.Range('E44').Value = '  +0.02%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '19.66'
$ws.Range('E45').Value = '  +0.71%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0510'
$ws.Range('E46').Value = '  -1.23%  '
$ws.Range('E47').Value = '  -1.70%  '
$ws.Range('E48').Value = '  -2.72%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '11.05'
$ws.Range('E49').Value = '  -0.12%  '
$ws.Range('E50').Value = '  -0.04%  '
$ws.Range('E51').Value = '  -2.79%  '
